# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 589
$ws1.Range("F4").Value = 1269
$ws1.Range("F5").Value = 1113
$ws1.Range("F6").Value = 14128
$ws1.Range("F7").Value = 15726
$ws1.Range("F8").Value = 10
$ws1.Range("F9").Value = 56
$ws1.Range("F11").Value = 192
$ws1.Range("F18").Value = 81
$ws1.Range("F19").Value = 30
$ws1.Range("F20").Value = 1222
$ws1.Range("F22").Value = 66
$ws1.Range("F23").Value = 12
$ws1.Range("F24").Value = 6197
$ws1.Range("F26").Value = 1095
$ws1.Range("F27").Value = 5573
$ws1.Range("F28").Value = 75
$ws1.Range("F30").Value = 123
$ws1.Range("F31").Value = 4547
$ws1.Range("F32").Value = 3

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 589
$ws4.Range("F4").Value = 1269
$ws4.Range("F5").Value = 1113
$ws4.Range("F6").Value = 14128
$ws4.Range("F7").Value = 15727
$ws4.Range("F8").Value = 10
$ws4.Range("F9").Value = 56
$ws4.Range("F11").Value = 192
$ws4.Range("F18").Value = 81
$ws4.Range("F19").Value = 30
$ws4.Range("F20").Value = 1222
$ws4.Range("F22").Value = 66
$ws4.Range("F24").Value = 12
$ws4.Range("F25").Value = 6197
$ws4.Range("F27").Value = 1095
$ws4.Range("F28").Value = 5573
$ws4.Range("F29").Value = 75
$ws4.Range("F31").Value = 123
$ws4.Range("F32").Value = 4547
$ws4.Range("F33").Value = 3

$wb.Save()
